$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Add a new row to the "testSimple" table (mirrors the existing step1/run row)
$ws.Range("B16").Value = "step2"
$ws.Range("C16").Value = "run2"

# Match the formatting of the row above (step1/run)
$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
